$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.004.20'
$ws.Range("E2").Value = '  -6.17%  '
$ws.Range("D3").Value = '3.439.29'
$ws.Range("E3").Value = '  -7.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.23'
$ws.Range("E5").Value = '  -10.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.24'
$ws.Range("E6").Value = '  -6.38%  '
$ws.Range("D7").Value = '3.441.26'
$ws.Range("E7").Value = '  -7.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  -6.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.636'
$ws.Range("E10").Value = '  -12.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.139'
$ws.Range("E11").Value = '  -14.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.58'
$ws.Range("E12").Value = '  -16.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  -15.00%  '
$ws.Range("E14").Value = '  -13.17%  '
$ws.Range("D15").Value = '3.986.29'
$ws.Range("E15").Value = '  -7.60%  '
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").Value = '3.431.96'
$ws.Range("E17").Value = '  -7.78%  '
$ws.Range("D18").Value = '64.670.86'
$ws.Range("E18").Value = '  -6.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.54'
$ws.Range("E19").Value = '  -9.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.51'
$ws.Range("E20").Value = '  -10.98%  '
$ws.Range("E21").Value = '  -11.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.41'
$ws.Range("E22").Value = '  -9.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.02'
$ws.Range("E23").Value = '  -12.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.68'
$ws.Range("E24").Value = '  -8.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.57'
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.97'
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("E27").Value = '  -10.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.66'
$ws.Range("E28").Value = '  -9.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.38'
$ws.Range("E29").Value = '  -11.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.39'
$ws.Range("E30").Value = '  -13.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '29.77'
$ws.Range("E31").Value = '  -10.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.97'
$ws.Range("E32").Value = '  -10.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '598.68'
$ws.Range("E33").Value = '  -5.37%  '
$ws.Range("E34").Value = '  -8.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '62.34'
$ws.Range("E35").Value = '  -5.33%  '
$ws.Range("E36").Value = '  -13.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.68'
$ws.Range("E37").Value = '  -13.67%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.384'
$ws.Range("E39").Value = '  -7.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").Value = '0.0₃0700'
$ws.Range("E41").Value = '  -15.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.126'
$ws.Range("E42").Value = '  -10.87%  '
$ws.Range("D43").Value = '2.883.72'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -12.74%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  -8.93%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.06'
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0386'
$ws.Range("E47").Value = '  -13.85%  '
$ws.Range("E48").Value = '  -10.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.69'
$ws.Range("E49").Value = '  -3.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.44'
$ws.Range("E50").Value = '  -11.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.98'
$ws.Range("E51").Value = '  -13.09%  '
